$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the new columns (AD, AE, AF), cloning the formatting of the
# existing header cell AC1 (bold, bordered, centered header style) via copy/paste.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for every data row (rows 2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
